$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Interface Cable SMA to U.FL / WRL-09145 ---
$ws.Range("A10").Value = "Interface Cable SMA to U.FL"
$ws.Range("B10").Value = "WRL-09145"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 4.95

# --- Row 11: GPS/GNSS Magnetic Mount Antenna SMA - 3m / GPS-14986 ---
$ws.Range("A11").Value = "GPS/GNSS Magnetic Mount Antenna SMA - 3m"
$ws.Range("B11").Value = "GPS-14986"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 12.95

# --- Row 12: GPS Antenna Ground Plate / DIY ---
$ws.Range("A12").Value = "GPS Antenna Ground Plate"
$ws.Range("B12").Value = "DIY"

# --- Row 13: Mast to get Base station antenna above houses ---
$ws.Range("A13").Value = "Mast to get Base station antenna above houses"

# Add hyperlinks for the new part numbers (B10, B11, B12), then restore the
# "Hyperlink" cell style (Add() creates its own explicit style record).
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.sparkfun.com/products/09145")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.sparkfun.com/products/14986")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.sparkfun.com/products/diy")
$ws.Range("B10").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"

# Column A needs to widen to fit the new, longer descriptions.
$ws.Columns("A:A").AutoFit()

# Move the active selection, mirroring the author's final cursor position.
$ws.Range("A18").Select()
